# Refresh the EC (account-statement) database and add part 1 of the new
# period data: the "Periodo Mora" list is now stored in ascending order
# (1701 .. 2003) instead of descending, and the matching "Valor Mora" /
# "Salario Basico" amounts for each period row are updated.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "1701"
$ws.Range("F16").Value = 27600
$ws.Range("G16").Value = 781242
$ws.Range("E17").Value = "1702"
$ws.Range("F17").Value = 27600
$ws.Range("G17").Value = 781242
$ws.Range("E18").Value = "1703"
$ws.Range("F18").Value = 27600
$ws.Range("G18").Value = 781242
$ws.Range("E19").Value = "1704"
$ws.Range("F19").Value = 27600
$ws.Range("G19").Value = 781242
$ws.Range("E20").Value = "1705"
$ws.Range("F20").Value = 27600
$ws.Range("G20").Value = 781242
$ws.Range("E21").Value = "1706"
$ws.Range("F21").Value = 27600
$ws.Range("G21").Value = 781242
$ws.Range("E22").Value = "1707"
$ws.Range("F22").Value = 27600
$ws.Range("G22").Value = 781242
$ws.Range("E23").Value = "1708"
$ws.Range("F23").Value = 27600
$ws.Range("G23").Value = 781242
$ws.Range("E24").Value = "1709"
$ws.Range("F24").Value = 27600
$ws.Range("G24").Value = 781242
$ws.Range("E25").Value = "1710"
$ws.Range("F25").Value = 27600
$ws.Range("G25").Value = 781242
$ws.Range("E26").Value = "1711"
$ws.Range("F26").Value = 27600
$ws.Range("G26").Value = 781242
$ws.Range("E27").Value = "1712"
$ws.Range("F27").Value = 27600
$ws.Range("G27").Value = 781242
$ws.Range("E28").Value = "1801"
$ws.Range("F28").Value = 27600
$ws.Range("G28").Value = 781242
$ws.Range("E29").Value = "1802"
$ws.Range("F29").Value = 27600
$ws.Range("G29").Value = 781242
$ws.Range("E30").Value = "1803"
$ws.Range("F30").Value = 27600
$ws.Range("G30").Value = 781242
$ws.Range("E31").Value = "1804"
$ws.Range("F31").Value = 27600
$ws.Range("G31").Value = 781242
$ws.Range("E32").Value = "1805"
$ws.Range("F32").Value = 27600
$ws.Range("G32").Value = 781242
$ws.Range("E33").Value = "1806"
$ws.Range("F33").Value = 27600
$ws.Range("G33").Value = 781242
$ws.Range("E34").Value = "1807"
$ws.Range("F34").Value = 27600
$ws.Range("G34").Value = 781242
$ws.Range("E35").Value = "1808"
$ws.Range("F35").Value = 27600
$ws.Range("G35").Value = 781242
$ws.Range("E36").Value = "1809"
$ws.Range("F36").Value = 31249
$ws.Range("G36").Value = 781242
$ws.Range("E37").Value = "1810"
$ws.Range("F37").Value = 31249
$ws.Range("G37").Value = 781242
$ws.Range("E38").Value = "1811"
$ws.Range("F38").Value = 31249
$ws.Range("G38").Value = 781242
$ws.Range("E39").Value = "1812"
$ws.Range("F39").Value = 31249
$ws.Range("G39").Value = 781242
$ws.Range("E40").Value = "1901"
$ws.Range("F40").Value = 31249
$ws.Range("G40").Value = 781242
$ws.Range("E41").Value = "1902"
$ws.Range("F41").Value = 31249
$ws.Range("G41").Value = 781242
$ws.Range("E42").Value = "1903"
$ws.Range("F42").Value = 31249
$ws.Range("G42").Value = 781242
$ws.Range("E43").Value = "1904"
$ws.Range("F43").Value = 31249
$ws.Range("G43").Value = 781242
$ws.Range("E44").Value = "1905"
$ws.Range("F44").Value = 31249
$ws.Range("G44").Value = 781242
$ws.Range("E45").Value = "1906"
$ws.Range("F45").Value = 31249
$ws.Range("G45").Value = 781242
$ws.Range("E46").Value = "1907"
$ws.Range("F46").Value = 31249
$ws.Range("G46").Value = 781242
$ws.Range("E47").Value = "1908"
$ws.Range("F47").Value = 31249
$ws.Range("G47").Value = 781242
$ws.Range("E48").Value = "1909"
$ws.Range("F48").Value = 31249
$ws.Range("G48").Value = 781242
$ws.Range("E49").Value = "1910"
$ws.Range("F49").Value = 31249
$ws.Range("G49").Value = 781242
$ws.Range("E50").Value = "1911"
$ws.Range("F50").Value = 31249
$ws.Range("G50").Value = 781242
$ws.Range("E51").Value = "1912"
$ws.Range("F51").Value = 31249
$ws.Range("G51").Value = 781242
$ws.Range("E52").Value = "2001"
$ws.Range("F52").Value = 31249
$ws.Range("G52").Value = 781242
$ws.Range("E53").Value = "2002"
$ws.Range("F53").Value = 31249
$ws.Range("G53").Value = 781242
$ws.Range("E54").Value = "2003"
$ws.Range("F54").Value = 31249
$ws.Range("G54").Value = 781242

# Nudge the logo image left to line up with the narrower data columns
# (net size/shape stays the same, only its left position moves).
$shp = $ws.Shapes.Item(1)
$shp.Left = 59.09055118110236
$shp.Width = 76.81889763779527
$shp.Height = 48.188976377952756
